$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row labels: *_old -> *_FV2410, *_new -> *_FV2504 ---
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_old$", "_FV2410")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cur = $cell.Value()
    $cell.Value = ($cur -replace "_new$", "_FV2504")
}

# --- Turn the data range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
